$wb = $excel.ActiveWorkbook

# --- Login Page sheet: insert a new "Business Usecase ID" column before the
#     existing "is Automated?" column, shifting Module/Title/Steps/etc right.
$wsLogin = $wb.Worksheets.Item(2)
$wsLogin.Columns("B:B").Insert()
$wsLogin.Range("B2").Value = "Business Usecase ID"

# --- Rename "Test Plan" -> "Business use cases" and populate it with the new
#     business-use-case table.
$wsBiz = $wb.Worksheets.Item(1)
$wsBiz.Name = "Business use cases"

$wsBiz.Range("A1").Value = "ID"
$wsBiz.Range("B1").Value = "Use case"
$wsBiz.Range("C1").Value = "Acceptance criteria"

$wsBiz.Range("B2").Value = "As a shopper, I want to login to the Swag Labs site, so I can view items and place orders."
$wsBiz.Range("B3").Value = "As a shopper, I want to select and remove items, so I can upadate the cart with the items I wish to buy"
$wsBiz.Range("B4").Value = "As a shopper, I want to checkout the items in my shopping cart, so I can purchase the items selected"

$wsBiz.Columns("B:B").ColumnWidth = 47.6640625
$wsBiz.Columns("C:C").ColumnWidth = 42

$wsBiz.Range("B1").WrapText = $true
$wsBiz.Range("B3").WrapText = $true
$wsBiz.Range("B4").WrapText = $true

$wsBiz.Range("B2:C2").WrapText = $true
$wsBiz.Range("B2:C2").VerticalAlignment = -4108

$wsBiz.Rows("1:1").RowHeight = 17
$wsBiz.Rows("2:4").RowHeight = 34

$wsBiz.Range("D1").Select()
$wsLogin.Range("E54").Select()
